$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.214.91"
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").Value = "2.055.51"
$ws.Range("E3").Value = "  -0.79%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.82"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.47%  "

$ws.Range("E6").Value = "  -1.30%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.53"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -3.82%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.386"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.30%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0786"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.39%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.108"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.20%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.87"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.72%  "

$ws.Range("D13").Value = "2.353.60"
$ws.Range("E13").Value = "  -0.89%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.836"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.04%  "

$ws.Range("E15").Value = "  +5.48%  "

$ws.Range("D16").Value = "2.054.41"
$ws.Range("E16").Value = "  -0.77%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.00"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +18.31%  "

$ws.Range("D18").Value = "37.147.14"
$ws.Range("E18").Value = "  -0.10%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "75.19"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.58%  "

$ws.Range("D20").Value = "0.0₃0901"
$ws.Range("E20").Value = "  -3.13%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.39"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.52%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "237.60"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.89%  "

$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.49"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.20%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.20"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.63%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.67"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.67%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.41"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.37%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.17"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.05%  "

$ws.Range("E29").Value = "  -0.91%  "

$ws.Range("E30").Value = "  +1.41%  "

$ws.Range("E31").Value = "  +3.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0622"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.20%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.55"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.13%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0898"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.06%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.29"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.27%  "

$ws.Range("E37").Value = "  -0.12%  "

$ws.Range("E38").Value = "  -1.46%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.104"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.45%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.16"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +12.00%  "

$ws.Range("E41").Value = "  +10.45%  "

$ws.Range("E42").Value = "  -1.60%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.32"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.37%  "

$ws.Range("E44").Value = "  -0.93%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "96.35"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.23%  "

$ws.Range("E46").Value = "  -1.99%  "

$ws.Range("E47").Value = "  -1.37%  "

$ws.Range("D48").Value = "1.284.44"
$ws.Range("E48").Value = "  -1.59%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.85"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.18%  "

$ws.Range("D50").Value = "2.241.35"
$ws.Range("E50").Value = "  -2.16%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.59"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -19.48%  "
